# 陳碧涵 財產申報表 — sheet 5 (保險/insurance) and sheet 6 (債務/debt)
# Add the standard trailer columns (property_category/category/date/
# legislator_name/legislator_id/source_file/index) and fix up the header
# rows to use proper field-name labels, per commit "#5: insurance, claim,
# debt, investment done".
#
# Note: the literal text "2013-12-31" looks like a date to Excel's value
# parser, so writing it straight into .Value gets auto-converted to a date
# serial. To keep it as plain text (matching the original data, which is
# a shared string, not a date), we write it as a text formula and then
# collapse the formula down to its literal text result via PasteSpecial
# (values-only) — this avoids both the date coercion AND any stray
# number-format/style being attached to the cell.

$wb = $excel.ActiveWorkbook

function Set-TextDate($range) {
    $range.Formula = '="2013-12-31"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# Sheet 5: 保險 (insurance)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Header row (field names)
$ws5.Range("B1").Value = "company"
$ws5.Range("C1").Value = "name"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "property_category"
$ws5.Range("F1").Value = "category"
$ws5.Range("G1").Value = "date"
$ws5.Range("H1").Value = "legislator_name"
$ws5.Range("I1").Value = "legislator_id"
$ws5.Range("J1").Value = "source_file"
$ws5.Range("K1").Value = "index"

# Row 2 (index 131)
$ws5.Range("B2").Value = "南山人壽"
$ws5.Range("C2").Value = "南山金福利21年期養老壽險"
$ws5.Range("D2").Value = "陳碧涵"
$ws5.Range("E2").Value = "insurance"
$ws5.Range("F2").Value = "normal"
Set-TextDate $ws5.Range("G2")
$ws5.Range("H2").Value = "陳碧涵"
$ws5.Range("I2").Value = 1752
$ws5.Range("J2").Value = "tmp11ae1"
$ws5.Range("K2").Value = 131

# Row 3 (index 132)
$ws5.Range("B3").Value = "台銀人壽股份有限公司"
$ws5.Range("C3").Value = "鴻福還本終身壽險"
$ws5.Range("D3").Value = "陳碧涵"
$ws5.Range("E3").Value = "insurance"
$ws5.Range("F3").Value = "normal"
Set-TextDate $ws5.Range("G3")
$ws5.Range("H3").Value = "陳碧涵"
$ws5.Range("I3").Value = 1752
$ws5.Range("J3").Value = "tmp11ae1"
$ws5.Range("K3").Value = 132

# Row 4 (index 134)
$ws5.Range("B4").Value = "南山人壽"
$ws5.Range("C4").Value = "南山新康祥終身壽險"
$ws5.Range("D4").Value = "廖〇陽"
$ws5.Range("E4").Value = "insurance"
$ws5.Range("F4").Value = "normal"
Set-TextDate $ws5.Range("G4")
$ws5.Range("H4").Value = "陳碧涵"
$ws5.Range("I4").Value = 1752
$ws5.Range("J4").Value = "tmp11ae1"
$ws5.Range("K4").Value = 134

# ---------------------------------------------------------------------
# Sheet 6: 債務 (debt)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Header row (field names)
$ws6.Range("B1").Value = "species"
$ws6.Range("C1").Value = "debtor"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"
$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

# Row 2 (index 144)
$ws6.Range("H2").Value = "debt"
$ws6.Range("I2").Value = "normal"
Set-TextDate $ws6.Range("J2")
$ws6.Range("K2").Value = "陳碧涵"
$ws6.Range("L2").Value = 1752
$ws6.Range("M2").Value = "tmp11ae1"
$ws6.Range("N2").Value = 144

# Row 3 (index 145)
$ws6.Range("H3").Value = "debt"
$ws6.Range("I3").Value = "normal"
Set-TextDate $ws6.Range("J3")
$ws6.Range("K3").Value = "陳碧涵"
$ws6.Range("L3").Value = 1752
$ws6.Range("M3").Value = "tmp11ae1"
$ws6.Range("N3").Value = 145

# Row 4 (index 146)
$ws6.Range("H4").Value = "debt"
$ws6.Range("I4").Value = "normal"
Set-TextDate $ws6.Range("J4")
$ws6.Range("K4").Value = "陳碧涵"
$ws6.Range("L4").Value = 1752
$ws6.Range("M4").Value = "tmp11ae1"
$ws6.Range("N4").Value = 146

# Row 5 (index 147)
$ws6.Range("H5").Value = "debt"
$ws6.Range("I5").Value = "normal"
Set-TextDate $ws6.Range("J5")
$ws6.Range("K5").Value = "陳碧涵"
$ws6.Range("L5").Value = 1752
$ws6.Range("M5").Value = "tmp11ae1"
$ws6.Range("N5").Value = 147
